# Fruta / hortaliza, semanal
#
# Three new weekly price rows (dated 44488) are inserted at the top of the
# data block for this sheet (rows 131-133), pushing the existing rows
# 131-205 down to 134-208. The new rows follow the same constant columns
# (A,B,C,E,F,G,H,I,J) as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 131, shifting the rest of the
# data block (rows 131-205) down to rows 134-208.
$ws.Rows("131:133").Insert()

# --- New row 131 ---
$ws.Cells.Item(131, 1).Value = 5
$ws.Cells.Item(131, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(131, 3).Value = "Maule"
$ws.Cells.Item(131, 4).Value = 44488
$ws.Cells.Item(131, 5).Value = 7
$ws.Cells.Item(131, 6).Value = "Fruta"
$ws.Cells.Item(131, 7).Value = 100102
$ws.Cells.Item(131, 8).Value = "Cítricos"
$ws.Cells.Item(131, 9).Value = 100102004
$ws.Cells.Item(131, 10).Value = "Mandarina"
$ws.Cells.Item(131, 11).Value = "Murcott"
$ws.Cells.Item(131, 12).Value = "Primera"
$ws.Cells.Item(131, 13).Value = 10
$ws.Cells.Item(131, 14).Value = 140000
$ws.Cells.Item(131, 15).Value = 140000
$ws.Cells.Item(131, 16).Value = 140000
$ws.Cells.Item(131, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(131, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(131, 19).Value = 311
$ws.Cells.Item(131, 20).Value = 450

# --- New row 132 ---
$ws.Cells.Item(132, 1).Value = 5
$ws.Cells.Item(132, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(132, 3).Value = "Maule"
$ws.Cells.Item(132, 4).Value = 44488
$ws.Cells.Item(132, 5).Value = 7
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100102
$ws.Cells.Item(132, 8).Value = "Cítricos"
$ws.Cells.Item(132, 9).Value = 100102004
$ws.Cells.Item(132, 10).Value = "Mandarina"
$ws.Cells.Item(132, 11).Value = "Murcott"
$ws.Cells.Item(132, 12).Value = "Primera"
$ws.Cells.Item(132, 13).Value = 200
$ws.Cells.Item(132, 14).Value = 6000
$ws.Cells.Item(132, 15).Value = 6000
$ws.Cells.Item(132, 16).Value = 6000
$ws.Cells.Item(132, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(132, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(132, 19).Value = 333
$ws.Cells.Item(132, 20).Value = 18

# --- New row 133 ---
$ws.Cells.Item(133, 1).Value = 5
$ws.Cells.Item(133, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(133, 3).Value = "Maule"
$ws.Cells.Item(133, 4).Value = 44488
$ws.Cells.Item(133, 5).Value = 7
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100102
$ws.Cells.Item(133, 8).Value = "Cítricos"
$ws.Cells.Item(133, 9).Value = 100102004
$ws.Cells.Item(133, 10).Value = "Mandarina"
$ws.Cells.Item(133, 11).Value = "Murcott"
$ws.Cells.Item(133, 12).Value = "Segunda"
$ws.Cells.Item(133, 13).Value = 160
$ws.Cells.Item(133, 14).Value = 4000
$ws.Cells.Item(133, 15).Value = 4000
$ws.Cells.Item(133, 16).Value = 4000
$ws.Cells.Item(133, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(133, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(133, 19).Value = 222
$ws.Cells.Item(133, 20).Value = 18
